$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 580 ("2026/12/29" block),
# shifting all subsequent rows (580-621) down by one, to 581-622.
$ws.Rows.Item(580).Insert()

# Populate the newly inserted row 580 with the new daily entry
# (date 2026/01/08, Thursday, hour 13, ranking 20).
$dateCell = $ws.Range("A580")
# Force text storage for the date-like string so Excel does not
# auto-convert it into a date serial number (column A stores plain
# text dates elsewhere in the sheet), then restore the default
# "Normal" style so no stray number-format style is left behind.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/01/08"
$dateCell.Style = "Normal"

$ws.Range("B580").Value = "木"
$ws.Range("C580").Value = 13
$ws.Range("D580").Value = 20
